# Apply updated crypto market data (prices and 1h volume %) per GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "69.458.05"
$ws.Range("E2").Value = "  -2.48%  "

# Row 3
$ws.Range("D3").Value = "3.697.64"
$ws.Range("E3").Value = "  -3.20%  "

# Row 4
$ws.Range("E4").Value = "  -0.05%  "

# Row 5
$ws.Range("D5").Value = "'694.21"
$ws.Range("E5").Value = "  -1.51%  "

# Row 6
$ws.Range("D6").Value = "'163.54"
$ws.Range("E6").Value = "  -5.05%  "

# Row 7
$ws.Range("D7").Value = "3.696.16"
$ws.Range("E7").Value = "  -3.21%  "

# Row 8
$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  -0.17%  "

# Row 9
$ws.Range("E9").Value = "  -4.28%  "

# Row 10
$ws.Range("E10").Value = "  -7.91%  "

# Row 11
$ws.Range("D11").Value = "'7.40"
$ws.Range("E11").Value = "  -1.93%  "

# Row 12
$ws.Range("E12").Value = "  -3.78%  "

# Row 13
$ws.Range("D13").Value = "'0.0000241"
$ws.Range("E13").Value = "  -4.74%  "

# Row 14
$ws.Range("D14").Value = "'33.54"
$ws.Range("E14").Value = "  -6.86%  "

# Row 15
$ws.Range("D15").Value = "4.318.84"
$ws.Range("E15").Value = "  -3.26%  "

# Row 16
$ws.Range("D16").Value = "3.696.04"
$ws.Range("E16").Value = "  -3.12%  "

# Row 17
$ws.Range("D17").Value = "69.519.39"
$ws.Range("E17").Value = "  -2.38%  "

# Row 18
$ws.Range("E18").Value = "  -0.87%  "

# Row 19
$ws.Range("D19").Value = "'16.32"
$ws.Range("E19").Value = "  -7.00%  "

# Row 20
$ws.Range("E20").Value = "  -7.57%  "

# Row 21
$ws.Range("D21").Value = "'483.14"
$ws.Range("E21").Value = "  -5.56%  "

# Row 22
$ws.Range("D22").Value = "'10.01"
$ws.Range("E22").Value = "  -6.82%  "

# Row 23
$ws.Range("E23").Value = "  -7.61%  "

# Row 24
$ws.Range("D24").Value = "'80.16"
$ws.Range("E24").Value = "  -4.82%  "

# Row 25
$ws.Range("D25").Value = "3.840.16"
$ws.Range("E25").Value = "  -3.31%  "

# Row 26
$ws.Range("E26").Value = "  -9.23%  "

# Row 27
$ws.Range("E27").Value = "  +0.01%  "

# Row 28
$ws.Range("E28").Value = "  -5.00%  "

# Row 29
$ws.Range("D29").Value = "'9.58"
$ws.Range("E29").Value = "  -8.02%  "

# Row 30
$ws.Range("E30").Value = "  -9.81%  "

# Row 31
$ws.Range("E31").Value = "  -9.87%  "

# Row 32
$ws.Range("D32").Value = "'6.89"
$ws.Range("E32").Value = "  -7.28%  "

# Row 33
$ws.Range("E33").Value = "  -7.11%  "

# Row 34
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").Value = "'27.15"
$ws.Range("E34").Value = "  -6.74%  "

# Row 35
$ws.Range("B35").Value = "Binance-PegBSC-USD"
$ws.Range("C35").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  +0.28%  "

# Row 36
$ws.Range("E36").Value = "  -4.12%  "

# Row 37
$ws.Range("D37").Value = "3.661.98"
$ws.Range("E37").Value = "  -3.18%  "

# Row 38
$ws.Range("D38").Value = "'8.53"
$ws.Range("E38").Value = "  -7.06%  "

# Row 39
$ws.Range("D39").Value = "'6.38"
$ws.Range("E39").Value = "  +6.08%  "

# Row 40
$ws.Range("E40").Value = "  -1.86%  "

# Row 41
$ws.Range("D41").Value = "'0.0937"
$ws.Range("E41").Value = "  -7.52%  "

# Row 42
$ws.Range("E42").Value = "  +0.00%  "

# Row 43
$ws.Range("E43").Value = "  -0.12%  "

# Row 44
$ws.Range("D44").Value = "'0.954"
$ws.Range("E44").Value = "  -6.83%  "

# Row 45
$ws.Range("D45").Value = "'163.92"
$ws.Range("E45").Value = "  -4.28%  "

# Row 46
$ws.Range("D46").Value = "'48.06"
$ws.Range("E46").Value = "  -3.11%  "

# Row 47
$ws.Range("B47").Value = "dogwifhat"
$ws.Range("C47").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D47").Value = "'2.84"
$ws.Range("E47").Value = "  -14.25%  "

# Row 48
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").Value = "'30.13"
$ws.Range("E48").Value = "  +2.31%  "

# Row 49
$ws.Range("B49").Value = "ONDO"
$ws.Range("C49").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D49").Value = "'1.36"
$ws.Range("E49").Value = "  +0.28%  "

# Row 50
$ws.Range("B50").Value = "SuiNetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D50").Value = "'1.15"
$ws.Range("E50").Value = "  +0.44%  "

# Row 51
$ws.Range("B51").Value = "FLOKI"
$ws.Range("C51").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D51").Value = "'0.000287"
$ws.Range("E51").Value = "  -7.39%  "

